# Apply updated odds values as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 update
$ws.Range("G3").Value = 2.1

# Row 6 updates
$ws.Range("G6").Value  = 3.1
$ws.Range("I6").Value  = 2.52
$ws.Range("J6").Value  = 3.8
$ws.Range("K6").Value  = 1.83
$ws.Range("L6").Value  = 3.2
$ws.Range("S6").Value  = 1.55
$ws.Range("T6").Value  = 2.15
$ws.Range("W6").Value  = 6.9
$ws.Range("Y6").Value  = 11.5
$ws.Range("Z6").Value  = 45
$ws.Range("AA6").Value = 35
$ws.Range("AB6").Value = 50
$ws.Range("AC6").Value = 5.9
$ws.Range("AH6").Value = 6.3
$ws.Range("AI6").Value = 11.5
$ws.Range("AJ6").Value = 10
$ws.Range("AK6").Value = 29
$ws.Range("AL6").Value = 25
$ws.Range("AM6").Value = 40
$ws.Range("AN6").Value = 4.75
$ws.Range("AO6").Value = 18.5
$ws.Range("AP6").Value = 29
$ws.Range("AQ6").Value = 100
$ws.Range("AR6").Value = 175
$ws.Range("AS6").Value = 450
$ws.Range("AV6").Value = 90
$ws.Range("AW6").Value = 4.2
$ws.Range("AX6").Value = 14
$ws.Range("AY6").Value = 24
$ws.Range("AZ6").Value = 65
$ws.Range("BA6").Value = 110
$ws.Range("BB6").Value = 350

$wb.Save()
